$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Cells.Item(2,4)
$c.NumberFormat = "@"
$c.Value = '301.69'
$c.Style = "Normal"
$c = $ws.Cells.Item(2,5)
$c.NumberFormat = "@"
$c.Value = '-0.85%'
$c.Style = "Normal"

# Row 3
$c = $ws.Cells.Item(3,4)
$c.NumberFormat = "@"
$c.Value = '31.48'
$c.Style = "Normal"
$c = $ws.Cells.Item(3,5)
$c.NumberFormat = "@"
$c.Value = '-3.30%'
$c.Style = "Normal"

# Row 4
$c = $ws.Cells.Item(4,4)
$c.NumberFormat = "@"
$c.Value = '5.149'
$c.Style = "Normal"
$c = $ws.Cells.Item(4,5)
$c.NumberFormat = "@"
$c.Value = '-2.77%'
$c.Style = "Normal"

# Row 5
$c = $ws.Cells.Item(5,4)
$c.NumberFormat = "@"
$c.Value = '0.07418'
$c.Style = "Normal"
$c = $ws.Cells.Item(5,5)
$c.NumberFormat = "@"
$c.Value = '-1.03%'
$c.Style = "Normal"

# Row 6
$c = $ws.Cells.Item(6,4)
$c.NumberFormat = "@"
$c.Value = '2.149'
$c.Style = "Normal"
$c = $ws.Cells.Item(6,5)
$c.NumberFormat = "@"
$c.Value = '42.59%'
$c.Style = "Normal"

# Row 7
$c = $ws.Cells.Item(7,4)
$c.NumberFormat = "@"
$c.Value = '7.923'
$c.Style = "Normal"
$c = $ws.Cells.Item(7,5)
$c.NumberFormat = "@"
$c.Value = '1.00%'
$c.Style = "Normal"

# Row 8
$c = $ws.Cells.Item(8,2)
$c.NumberFormat = "@"
$c.Value = 'MXToken'
$c.Style = "Normal"
$c = $ws.Cells.Item(8,3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c.Style = "Normal"
$c = $ws.Cells.Item(8,4)
$c.NumberFormat = "@"
$c.Value = '0.9261'
$c.Style = "Normal"
$c = $ws.Cells.Item(8,5)
$c.NumberFormat = "@"
$c.Value = '0.94%'
$c.Style = "Normal"

# Row 9
$c = $ws.Cells.Item(9,2)
$c.NumberFormat = "@"
$c.Value = 'WazirX'
$c.Style = "Normal"
$c = $ws.Cells.Item(9,3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$c.Style = "Normal"
$c = $ws.Cells.Item(9,4)
$c.NumberFormat = "@"
$c.Value = '0.1724'
$c.Style = "Normal"
$c = $ws.Cells.Item(9,5)
$c.NumberFormat = "@"
$c.Value = '1.70%'
$c.Style = "Normal"

# Row 10
$c = $ws.Cells.Item(10,2)
$c.NumberFormat = "@"
$c.Value = 'LiechtensteinCryptoassetsExchange'
$c.Style = "Normal"
$c = $ws.Cells.Item(10,3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$c.Style = "Normal"
$c = $ws.Cells.Item(10,4)
$c.NumberFormat = "@"
$c.Value = '0.07681'
$c.Style = "Normal"
$c = $ws.Cells.Item(10,5)
$c.NumberFormat = "@"
$c.Value = '-2.39%'
$c.Style = "Normal"

# Row 11
$c = $ws.Cells.Item(11,2)
$c.NumberFormat = "@"
$c.Value = 'MandalaExchangeToken'
$c.Style = "Normal"
$c = $ws.Cells.Item(11,3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$c.Style = "Normal"
$c = $ws.Cells.Item(11,4)
$c.NumberFormat = "@"
$c.Value = '0.08180'
$c.Style = "Normal"
$c = $ws.Cells.Item(11,5)
$c.NumberFormat = "@"
$c.Value = '1.34%'
$c.Style = "Normal"

# Row 12
$c = $ws.Cells.Item(12,2)
$c.NumberFormat = "@"
$c.Value = 'BitrueCoin'
$c.Style = "Normal"
$c = $ws.Cells.Item(12,3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$c.Style = "Normal"
$c = $ws.Cells.Item(12,4)
$c.NumberFormat = "@"
$c.Value = '0.03027'
$c.Style = "Normal"
$c = $ws.Cells.Item(12,5)
$c.NumberFormat = "@"
$c.Value = '0.23%'
$c.Style = "Normal"

# Row 13
$c = $ws.Cells.Item(13,2)
$c.NumberFormat = "@"
$c.Value = 'BitMartToken'
$c.Style = "Normal"
$c = $ws.Cells.Item(13,3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$c.Style = "Normal"
$c = $ws.Cells.Item(13,4)
$c.NumberFormat = "@"
$c.Value = '0.09931'
$c.Style = "Normal"
$c = $ws.Cells.Item(13,5)
$c.NumberFormat = "@"
$c.Value = '0.35%'
$c.Style = "Normal"

# Row 14
$c = $ws.Cells.Item(14,2)
$c.NumberFormat = "@"
$c.Value = 'BitForexToken'
$c.Style = "Normal"
$c = $ws.Cells.Item(14,3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$c.Style = "Normal"
$c = $ws.Cells.Item(14,4)
$c.NumberFormat = "@"
$c.Value = '0.001489'
$c.Style = "Normal"
$c = $ws.Cells.Item(14,5)
$c.NumberFormat = "@"
$c.Value = '0.10%'
$c.Style = "Normal"

# Row 15
$c = $ws.Cells.Item(15,2)
$c.NumberFormat = "@"
$c.Value = 'TigerCash'
$c.Style = "Normal"
$c = $ws.Cells.Item(15,3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$c.Style = "Normal"
$c = $ws.Cells.Item(15,4)
$c.NumberFormat = "@"
$c.Value = '0.006127'
$c.Style = "Normal"
$c = $ws.Cells.Item(15,5)
$c.NumberFormat = "@"
$c.Value = '-1.17%'
$c.Style = "Normal"

# Row 16
$c = $ws.Cells.Item(16,2)
$c.NumberFormat = "@"
$c.Value = 'LEO'
$c.Style = "Normal"
$c = $ws.Cells.Item(16,3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$c.Style = "Normal"
$c = $ws.Cells.Item(16,4)
$c.NumberFormat = "@"
$c.Value = '3.466'
$c.Style = "Normal"
$c = $ws.Cells.Item(16,5)
$c.NumberFormat = "@"
$c.Value = '-0.11%'
$c.Style = "Normal"

# Row 17
$c = $ws.Cells.Item(17,2)
$c.NumberFormat = "@"
$c.Value = 'GateToken'
$c.Style = "Normal"
$c = $ws.Cells.Item(17,3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$c.Style = "Normal"
$c = $ws.Cells.Item(17,4)
$c.NumberFormat = "@"
$c.Value = '3.762'
$c.Style = "Normal"
$c = $ws.Cells.Item(17,5)
$c.NumberFormat = "@"
$c.Value = '-1.16%'
$c.Style = "Normal"

# Row 18
$c = $ws.Cells.Item(18,4)
$c.NumberFormat = "@"
$c.Value = '2.231'
$c.Style = "Normal"
$c = $ws.Cells.Item(18,5)
$c.NumberFormat = "@"
$c.Value = '0.03%'
$c.Style = "Normal"

# Row 19
$c = $ws.Cells.Item(19,4)
$c.NumberFormat = "@"
$c.Value = '0.3249'
$c.Style = "Normal"
$c = $ws.Cells.Item(19,5)
$c.NumberFormat = "@"
$c.Value = '-2.37%'
$c.Style = "Normal"

# Row 20
$c = $ws.Cells.Item(20,4)
$c.NumberFormat = "@"
$c.Value = '0.1335'
$c.Style = "Normal"
$c = $ws.Cells.Item(20,5)
$c.NumberFormat = "@"
$c.Value = '0.01%'
$c.Style = "Normal"

# Row 21
$c = $ws.Cells.Item(21,4)
$c.NumberFormat = "@"
$c.Value = '4.648'
$c.Style = "Normal"
$c = $ws.Cells.Item(21,5)
$c.NumberFormat = "@"
$c.Value = '3.73%'
$c.Style = "Normal"

# Row 22
$c = $ws.Cells.Item(22,4)
$c.NumberFormat = "@"
$c.Value = '0.04641'
$c.Style = "Normal"
$c = $ws.Cells.Item(22,5)
$c.NumberFormat = "@"
$c.Value = '0.66%'
$c.Style = "Normal"

# Row 24
$c = $ws.Cells.Item(24,4)
$c.NumberFormat = "@"
$c.Value = '0.001221'
$c.Style = "Normal"
$c = $ws.Cells.Item(24,5)
$c.NumberFormat = "@"
$c.Value = '0.08%'
$c.Style = "Normal"

# Row 25
$c = $ws.Cells.Item(25,4)
$c.NumberFormat = "@"
$c.Value = '0.004496'
$c.Style = "Normal"
$c = $ws.Cells.Item(25,5)
$c.NumberFormat = "@"
$c.Value = '0.94%'
$c.Style = "Normal"

# Row 26
$c = $ws.Cells.Item(26,4)
$c.NumberFormat = "@"
$c.Value = '0.0001298'
$c.Style = "Normal"
$c = $ws.Cells.Item(26,5)
$c.NumberFormat = "@"
$c.Value = '-7.13%'
$c.Style = "Normal"

# Row 27
$c = $ws.Cells.Item(27,5)
$c.NumberFormat = "@"
$c.Value = '7.69%'
$c.Style = "Normal"

# Row 39
$c = $ws.Cells.Item(39,4)
$c.NumberFormat = "@"
$c.Value = '0.01740'
$c.Style = "Normal"
$c = $ws.Cells.Item(39,5)
$c.NumberFormat = "@"
$c.Value = '-1.75%'
$c.Style = "Normal"

# Row 40
$c = $ws.Cells.Item(40,4)
$c.NumberFormat = "@"
$c.Value = '0.04542'
$c.Style = "Normal"
$c = $ws.Cells.Item(40,5)
$c.NumberFormat = "@"
$c.Value = '-0.08%'
$c.Style = "Normal"

# Row 41
$c = $ws.Cells.Item(41,4)
$c.NumberFormat = "@"
$c.Value = '0.007128'
$c.Style = "Normal"
$c = $ws.Cells.Item(41,5)
$c.NumberFormat = "@"
$c.Value = '-0.64%'
$c.Style = "Normal"

# Row 42
$c = $ws.Cells.Item(42,4)
$c.NumberFormat = "@"
$c.Value = '0.1348'
$c.Style = "Normal"
$c = $ws.Cells.Item(42,5)
$c.NumberFormat = "@"
$c.Value = '0.21%'
$c.Style = "Normal"

# Row 43
$c = $ws.Cells.Item(43,4)
$c.NumberFormat = "@"
$c.Value = '0.002206'
$c.Style = "Normal"
$c = $ws.Cells.Item(43,5)
$c.NumberFormat = "@"
$c.Value = '0.02%'
$c.Style = "Normal"

# Row 44
$c = $ws.Cells.Item(44,4)
$c.NumberFormat = "@"
$c.Value = '0.01095'
$c.Style = "Normal"
$c = $ws.Cells.Item(44,5)
$c.NumberFormat = "@"
$c.Value = '-14.49%'
$c.Style = "Normal"

# Row 45
$c = $ws.Cells.Item(45,4)
$c.NumberFormat = "@"
$c.Value = '0.00006276'
$c.Style = "Normal"
$c = $ws.Cells.Item(45,5)
$c.NumberFormat = "@"
$c.Value = '1.04%'
$c.Style = "Normal"

# Row 46
$c = $ws.Cells.Item(46,5)
$c.NumberFormat = "@"
$c.Value = '-46.08%'
$c.Style = "Normal"

# Row 47
$c = $ws.Cells.Item(47,4)
$c.NumberFormat = "@"
$c.Value = '0.7383'
$c.Style = "Normal"
$c = $ws.Cells.Item(47,5)
$c.NumberFormat = "@"
$c.Value = '4.09%'
$c.Style = "Normal"
